# Horarios actualizados - Linea 141 - scrape refresh 05:57:38
# Updates Hora_Scrap / Minutos (and a few re-sequenced rows) across all three sheets

function Set-Row {
    param($ws, [int]$Row, $HoraScrap, $HoraLlegada, $Linea, $Minutos, $Parada)
    $ws.Cells.Item($Row, 1).Value = $HoraScrap
    $ws.Cells.Item($Row, 2).Value = $HoraLlegada
    $ws.Cells.Item($Row, 3).Value = $Linea
    $ws.Cells.Item($Row, 4).Value = $Minutos
    $ws.Cells.Item($Row, 5).Value = $Parada
}

$wb = $excel.ActiveWorkbook

# ----- Sheet: LP1912 -----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 05:57:38"
$ws.Range("A3").Value = "Total filas: 52"
$ws.Range("A8").Value = "04:44:46"
$ws.Range("C8").Value = "215_EL PELIGRO"
$ws.Range("D8").Value = 2
$ws.Range("C9").Value = "15_ABASTO"
$ws.Range("A10").Value = "03:52:04"
$ws.Range("C10").Value = "215A_EL PATO"
$ws.Range("D10").Value = 54
$ws.Range("A27").Value = "05:57:38"
$ws.Range("D27").Value = 3
$ws.Range("A29").Value = "05:57:38"
$ws.Range("D29").Value = 6
$ws.Range("A31").Value = "05:57:38"
$ws.Range("D31").Value = 13
$ws.Range("A34").Value = "05:57:38"
$ws.Range("D34").Value = 27
$ws.Range("A35").Value = "05:57:38"
$ws.Range("D35").Value = 30
$ws.Range("A38").Value = "05:57:38"
$ws.Range("D38").Value = 34
$ws.Range("A39").Value = "05:57:38"
$ws.Range("D39").Value = 34
$ws.Range("A41").Value = "05:57:38"
$ws.Range("D41").Value = 42
$ws.Range("A44").Value = "05:57:38"
$ws.Range("D44").Value = 57
$ws.Range("A45").Value = "05:57:38"
$ws.Range("D45").Value = 67
$ws.Range("A46").Value = "05:57:38"
$ws.Range("D46").Value = 69
$ws.Range("A48").Value = "05:57:38"
$ws.Range("D48").Value = 76
$ws.Range("A50").Value = "05:57:38"
$ws.Range("B50").Value = "07:20"
$ws.Range("D50").Value = 83
$ws.Range("B51").Value = "07:21"
$ws.Range("C51").Value = "215A_EL PATO"
$ws.Range("D51").Value = 98
$ws.Range("A52").Value = "05:57:38"
$ws.Range("B52").Value = "07:29"
$ws.Range("C52").Value = "14_ABASTO"
$ws.Range("D52").Value = 92
$ws.Range("A53").Value = "05:57:38"
$ws.Range("B53").Value = "07:33"
$ws.Range("C53").Value = "23_HERNANDEZ"
$ws.Range("D53").Value = 96
Set-Row $ws 54 "05:57:38" "07:36" "27_EL RETIRO" 99 "LP1912"
Set-Row $ws 55 "05:57:38" "07:36" "17X38_ROMERO" 99 "LP1912"
Set-Row $ws 56 "05:57:38" "07:43" "10_OLMOS" 106 "LP1912"
Set-Row $ws 57 "05:57:38" "07:49" "15_ABASTO" 112 "LP1912"

# ----- Sheet: LP1912-215 -----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 05:57:38"
$ws.Range("A3").Value = "Total filas: 11"
$ws.Range("A10").Value = "05:57:38"
$ws.Range("D10").Value = 13
$ws.Range("A13").Value = "05:57:38"
$ws.Range("D13").Value = 69
$ws.Range("A15").Value = "05:57:38"
$ws.Range("B15").Value = "07:20"
$ws.Range("D15").Value = 83
Set-Row $ws 16 "05:43:29" "07:21" "215A_EL PATO" 98 "LP1912"

# ----- Sheet: 6203-6173 -----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 05:57:38"
$ws.Range("A7").Value = "05:57:38"
$ws.Range("D7").Value = 90
